# Installation Guide.docx - update the notebook file name/date stamp from
# "2019Dec"/"2019JOct" to "2021Jan" in the three places it appears:
#   1. Open the notebook "Geom Alg Palette 2019Dec"  -> "...2021Jan"
#   2. Open the notebook GeomAlg2019JOct src         -> GeomAlg2021Jan src
#   3. GeomAlg2019Dec src - Source file for package  -> GeomAlg2021Jan src - ...

$d = $word.ActiveDocument

function Replace-ExactRange($needle, $replacement) {
    $full = $d.Content.Text
    $idx = $full.IndexOf($needle)
    if ($idx -ge 0) {
        $r = $d.Range($idx, $idx + $needle.Length)
        $r.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $replacement, 2) | Out-Null
    }
}

# 1) Quoted notebook name in the first instruction paragraph.
Replace-ExactRange "Geom Alg Palette 2019Dec" "Geom Alg Palette 2021Jan"

# 2) Notebook / source file name referenced later in the instructions.
Replace-ExactRange "GeomAlg2019JOct src" "GeomAlg2021Jan src"

# 3) Bold file-name heading in the "Files:" list.
Replace-ExactRange "GeomAlg2019Dec src" "GeomAlg2021Jan src"
